$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (update through 20/09/2021), appended after the
# existing last row (374, serial date 44448).
$data = @(
    @(44449, 0, 1, 33.71544167228591),
    @(44450, 0, 1, 33.71544167228591),
    @(44451, 0, 1, 33.71544167228591),
    @(44452, 0, 0, 0),
    @(44453, 0, 0, 0),
    @(44454, 0, 0, 0),
    @(44455, 0, 0, 0),
    @(44456, 0, 0, 0),
    @(44457, 0, 0, 0),
    @(44458, 0, 0, 0),
    @(44459, 0, 0, 0)
)

$startRow = 375
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Carry the date-cell style (column A) from the previous last row down to
# the newly appended rows, matching the existing formatting pattern.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
$excel.CutCopyMode = $false
